$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5927688107317977
$ws.Range("C2").Value = 0.08198816812695497
$ws.Range("E2").Value = 0.1623990570935394
$ws.Range("F2").Value = 2.651338463334199
$ws.Range("G2").Value = 1.590116522096409
$ws.Range("H2").Value = 1.441423647620056
$ws.Range("I2").Value = 1.513088120598908
$ws.Range("J2").Value = 0.1108423745061238
$ws.Range("K2").Value = 0.3440564424601291
$ws.Range("L2").Value = 0.3815664171981723
$ws.Range("M2").Value = 0.2263012922041483

$ws.Range("B3").Value = 0.5625585361359242
$ws.Range("C3").Value = 0.07966545906609213
$ws.Range("E3").Value = 0.1623416491463523
$ws.Range("F3").Value = 2.650935579588278
$ws.Range("G3").Value = 1.593830716240106
$ws.Range("H3").Value = 1.448008630242668
$ws.Range("I3").Value = 1.52021051249455
$ws.Range("J3").Value = 0.1103978555219349
$ws.Range("K3").Value = 0.3148555271466478
$ws.Range("L3").Value = 0.3777868916016587
$ws.Range("M3").Value = 0.2198238735868969

$ws.Range("B4").Value = 0.5442621926219999
$ws.Range("C4").Value = 0.07821460717375572
$ws.Range("E4").Value = 0.1623506468434712
$ws.Range("F4").Value = 2.651825687326593
$ws.Range("G4").Value = 1.596867967201533
$ws.Range("H4").Value = 1.452571055002608
$ws.Range("I4").Value = 1.525174154289118
$ws.Range("J4").Value = 0.110122141760078
$ws.Range("K4").Value = 0.2970412192054965
$ws.Range("L4").Value = 0.3756132723900691
$ws.Range("M4").Value = 0.2159456075824586

$ws.Range("B5").Value = 0.5368703563343047
$ws.Range("C5").Value = 0.07761713939631676
$ws.Range("E5").Value = 0.1623654786881996
$ws.Range("F5").Value = 2.652474726800619
$ws.Range("G5").Value = 1.598295920654948
$ws.Range("H5").Value = 1.45456093461442
$ws.Range("I5").Value = 1.527345413321882
$ws.Range("J5").Value = 0.110009099068014
$ws.Range("K5").Value = 0.2898110106710448
$ws.Range("L5").Value = 0.3747646154625315
$ws.Range("M5").Value = 0.2143901844867777

$ws.Range("B6").Value = 0.5356468315958125
$ws.Range("C6").Value = 0.07751755314645692
$ws.Range("E6").Value = 0.162368617006873
$ws.Range("F6").Value = 2.652599800326215
$ws.Range("G6").Value = 1.598544520818763
$ws.Range("H6").Value = 1.454899246554405
$ws.Range("I6").Value = 1.527714922066352
$ws.Range("J6").Value = 0.1099902872498202
$ws.Range("K6").Value = 0.2886122167379028
$ws.Range("L6").Value = 0.3746259421378042
$ws.Range("M6").Value = 0.2141334216631137

$ws.Range("B7").Value = 0.5441622436649993
$ws.Range("C7").Value = 0.07820657478530535
$ws.Range("E7").Value = 0.162350801611133
$ws.Range("F7").Value = 2.651833280853523
$ws.Range("G7").Value = 1.596886454819483
$ws.Range("H7").Value = 1.452597362080354
$ws.Range("I7").Value = 1.525202835163157
$ws.Range("J7").Value = 0.1101206199942961
$ws.Range("K7").Value = 0.2969435911784899
$ws.Range("L7").Value = 0.3756016766599402
$ws.Range("M7").Value = 0.2159245291932379

$ws.Range("B8").Value = 0.5823000989115883
$ws.Range("C8").Value = 0.08119240992969878
$ws.Range("E8").Value = 0.1623701026465056
$ws.Range("F8").Value = 2.650963629773202
$ws.Range("G8").Value = 1.591240095302112
$ws.Range("H8").Value = 1.443586440432668
$ws.Range("I8").Value = 1.515421403969675
$ws.Range("J8").Value = 0.1106896873998551
$ws.Range("K8").Value = 0.3339642361362394
$ws.Range("L8").Value = 0.3802328036805704
$ws.Range("M8").Value = 0.2240474301395778

$ws.Range("B9").Value = 0.6590780939665422
$ws.Range("C9").Value = 0.08685290076903129
$ws.Range("E9").Value = 0.1627574207846152
$ws.Range("F9").Value = 2.658273901206243
$ws.Range("G9").Value = 1.586174778665679
$ws.Range("H9").Value = 1.430032212654837
$ws.Range("I9").Value = 1.500923155677903
$ws.Range("J9").Value = 0.1117831458663083
$ws.Range("K9").Value = 0.4074659662514364
$ws.Range("L9").Value = 0.3904756700188017
$ws.Range("M9").Value = 0.2407562690972895

$ws.Range("B10").Value = 0.7166838497737729
$ws.Range("C10").Value = 0.09089499802240653
$ws.Range("E10").Value = 0.1632530325303634
$ws.Range("F10").Value = 2.669131957960388
$ws.Range("G10").Value = 1.58612148934921
$ws.Range("H10").Value = 1.422579090505906
$ws.Range("I10").Value = 1.493124455442761
$ws.Range("J10").Value = 0.1125722828817466
$ws.Range("K10").Value = 0.4620120918519035
$ws.Range("L10").Value = 0.3987030562539502
$ws.Range("M10").Value = 0.2535027149940987

$ws.Range("B11").Value = 0.743147141038861
$ws.Range("C11").Value = 0.09270897652548626
$ws.Range("E11").Value = 0.1635239155933768
$ws.Range("F11").Value = 2.675261442618492
$ws.Range("G11").Value = 1.586895207426579
$ws.Range("H11").Value = 1.419731635768514
$ws.Range("I11").Value = 1.490195862413955
$ws.Range("J11").Value = 0.1129280879384069
$ws.Range("K11").Value = 0.486943669054142
$ws.Range("L11").Value = 0.4025971451451085
$ws.Range("M11").Value = 0.2594026119793256

$ws.Range("B12").Value = 0.7532048341371365
$ws.Range("C12").Value = 0.09339234734861179
$ws.Range("E12").Value = 0.1636329874008453
$ws.Range("F12").Value = 2.677753442198039
$ws.Range("G12").Value = 1.587303024140638
$ws.Range("H12").Value = 1.418731390251168
$ws.Range("I12").Value = 1.489175869390465
$ws.Range("J12").Value = 0.1130623549868517
$ws.Range("K12").Value = 0.4964013964464584
$ws.Range("L12").Value = 0.4040933847192179
$ws.Range("M12").Value = 0.2616512296944435

$ws.Range("B13").Value = 0.7510371092953392
$ws.Range("C13").Value = 0.09324532846724765
$ws.Range("E13").Value = 0.1636092086055214
$ws.Range("F13").Value = 2.677209148616498
$ws.Range("G13").Value = 1.587210085290948
$ws.Range("H13").Value = 1.418943341839068
$ws.Range("I13").Value = 1.489391585092349
$ws.Range("J13").Value = 0.1130334592024234
$ws.Range("K13").Value = 0.4943637689648597
$ws.Range("L13").Value = 0.4037701828623028
$ws.Range("M13").Value = 0.261166308935195

$ws.Range("B14").Value = 0.7439738624685504
$ws.Range("C14").Value = 0.09276526874704416
$ws.Range("E14").Value = 0.1635327590590947
$ws.Range("F14").Value = 2.675463037263398
$ws.Range("G14").Value = 1.586926456997773
$ws.Range("H14").Value = 1.419647781596481
$ws.Range("I14").Value = 1.490110163366566
$ws.Range("J14").Value = 0.1129391436275711
$ws.Range("K14").Value = 0.4877214295497367
$ws.Range("L14").Value = 0.4027198089245871
$ws.Range("M14").Value = 0.2595873182007011

$ws.Range("B15").Value = 0.7396521738688762
$ws.Range("C15").Value = 0.09247075748286449
$ws.Range("E15").Value = 0.1634867761235768
$ws.Range("F15").Value = 2.674415742371394
$ws.Range("G15").Value = 1.586767682669304
$ws.Range("H15").Value = 1.420089430517436
$ws.Range("I15").Value = 1.490561903336541
$ws.Range("J15").Value = 0.1128813112433171
$ws.Range("K15").Value = 0.4836549674417938
$ws.Range("L15").Value = 0.4020792376211801
$ws.Range("M15").Value = 0.2586220191595316

$ws.Range("B16").Value = 0.714959563021381
$ws.Range("C16").Value = 0.0907759532304695
$ws.Range("E16").Value = 0.1632362405326582
$ws.Range("F16").Value = 2.668755309225574
$ws.Range("G16").Value = 1.586086988957533
$ws.Range("H16").Value = 1.422776099924931
$ws.Range("I16").Value = 1.493328301483807
$ws.Range("J16").Value = 0.1125489652446063
$ws.Range("K16").Value = 0.4603851065674576
$ws.Range("L16").Value = 0.3984516031670182
$ws.Range("M16").Value = 0.2531191730438067

$ws.Range("B17").Value = 0.699877205111477
$ws.Range("C17").Value = 0.08972991126658059
$ws.Range("E17").Value = 0.1630941528358818
$ws.Range("F17").Value = 2.665587460649931
$ws.Range("G17").Value = 1.585873839628974
$ws.Range("H17").Value = 1.424563318976979
$ws.Range("I17").Value = 1.495183942300322
$ws.Range("J17").Value = 0.1123442595806701
$ws.Range("K17").Value = 0.4461398594335719
$ws.Range("L17").Value = 0.3962648468901051
$ws.Range("M17").Value = 0.2497692475348643

$ws.Range("B18").Value = 0.6912265593276743
$ws.Range("C18").Value = 0.08912592078551285
$ws.Range("E18").Value = 0.1630167058018372
$ws.Range("F18").Value = 2.663877430042319
$ws.Range("G18").Value = 1.585826344682957
$ws.Range("H18").Value = 1.42564239411675
$ws.Range("I18").Value = 1.496309529271656
$ws.Range("J18").Value = 0.1122262199224302
$ws.Range("K18").Value = 0.4379575217284071
$ws.Range("L18").Value = 0.3950213422459967
$ws.Range("M18").Value = 0.2478520184397723

$ws.Range("B19").Value = 0.6883017929177697
$ws.Range("C19").Value = 0.08892101832115173
$ws.Range("E19").Value = 0.1629912194876901
$ws.Range("F19").Value = 2.663317691942112
$ws.Range("G19").Value = 1.585823160245596
$ws.Range("H19").Value = 1.426016531144825
$ws.Range("I19").Value = 1.496700642495881
$ws.Range("J19").Value = 0.112186202835062
$ws.Range("K19").Value = 0.4351890480774045
$ws.Range("L19").Value = 0.3946027663314027
$ws.Range("M19").Value = 0.2472045236455003

$ws.Range("B20").Value = 0.7014802336899493
$ws.Range("C20").Value = 0.0898415057086126
$ws.Range("E20").Value = 0.1631088358278951
$ws.Range("F20").Value = 2.665913090357506
$ws.Range("G20").Value = 1.585888756279672
$ws.Range("H20").Value = 1.424367776783726
$ws.Range("I20").Value = 1.494980375437713
$ws.Range("J20").Value = 0.1123660818133168
$ws.Range("K20").Value = 0.4476551385285745
$ws.Range("L20").Value = 0.3964961560158002
$ws.Range("M20").Value = 0.2501248644521468

$ws.Range("B21").Value = 0.7460475194338585
$ws.Range("C21").Value = 0.09290636992053436
$ws.Range("E21").Value = 0.1635550381939268
$ws.Range("F21").Value = 2.675971276304466
$ws.Range("G21").Value = 1.587006648502182
$ws.Range("H21").Value = 1.419438753496181
$ws.Range("I21").Value = 1.489896684193546
$ws.Range("J21").Value = 0.1129668592084236
$ws.Range("K21").Value = 0.4896719955509354
$ws.Range("L21").Value = 0.4030277432211875
$ws.Range("M21").Value = 0.2600507148469049

$ws.Range("B22").Value = 0.7753880296432101
$ws.Range("C22").Value = 0.09488879185794019
$ws.Range("E22").Value = 0.1638844932600776
$ws.Range("F22").Value = 2.683540888700634
$ws.Range("G22").Value = 1.588406589591784
$ws.Range("H22").Value = 1.416672106829338
$ws.Range("I22").Value = 1.487092958178117
$ws.Range("J22").Value = 0.1133567652302183
$ws.Range("K22").Value = 0.5172295149215813
$ws.Range("L22").Value = 0.407422539909021
$ws.Range("M22").Value = 0.2666220220788205

$ws.Range("B23").Value = 0.7597091033229049
$ws.Range("C23").Value = 0.09383261775032281
$ws.Range("E23").Value = 0.1637052074247691
$ws.Range("F23").Value = 2.67940978821828
$ws.Range("G23").Value = 1.5875981439426
$ws.Range("H23").Value = 1.418107127790833
$ws.Range("I23").Value = 1.488541899666785
$ws.Range("J23").Value = 0.1131489192983786
$ws.Range("K23").Value = 0.5025127734906505
$ws.Range("L23").Value = 0.4050654688500543
$ws.Range("M23").Value = 0.2631071335083135

$ws.Range("B24").Value = 0.7007554408803287
$ws.Range("C24").Value = 0.08979106198508191
$ws.Range("E24").Value = 0.1631021844343081
$ws.Range("F24").Value = 2.665765526740245
$ws.Range("G24").Value = 1.585881778694585
$ws.Range("H24").Value = 1.424456020733004
$ws.Range("I24").Value = 1.495072224993557
$ws.Range("J24").Value = 0.1123562170755967
$ws.Range("K24").Value = 0.4469700576376567
$ws.Range("L24").Value = 0.3963915385057959
$ws.Range("M24").Value = 0.2499640629727082

$ws.Range("B25").Value = 0.6380962630893521
$ws.Range("C25").Value = 0.08534222344727738
$ws.Range("E25").Value = 0.1626154348707871
$ws.Range("F25").Value = 2.65533201995828
$ws.Range("G25").Value = 1.586901344460074
$ws.Range("H25").Value = 1.433258759159017
$ws.Range("I25").Value = 1.504344143160765
$ws.Range("J25").Value = 0.1114897975399565
$ws.Range("K25").Value = 0.3874857053892242
$ws.Range("L25").Value = 0.3875810035151517
$ws.Range("M25").Value = 0.2361531115194708

